# vce-youth.xlsx update: replace the old age-breakdown (14-18/19-21/22-29)
# chart data in columns V:Y with a single "WRC" program time series
# (2015-2019), and extend the blank chart-data placeholder rows (W/X)
# down through row 21 with a Comma number format, matching the new
# "family services over time" graph added to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: V1 "Age" -> "Program" --------------------------------------
$ws.Range("V1").Value = "Program"

# --- V2:V6 -> "WRC" (one program, five years of data) --------------------
$ws.Range("V2").Value = "WRC"
$ws.Range("V3").Value = "WRC"
$ws.Range("V4").Value = "WRC"
$ws.Range("V5").Value = "WRC"
$ws.Range("V6").Value = "WRC"

# --- W (Number) / X (Percent, now unused -> 0) / Y (Year, formula chain) -
$ws.Range("W2").Value = 259
$ws.Range("X2").Value = 0
# Y2 stays 2015 (unchanged)

$ws.Range("W3").Value = 250
$ws.Range("X3").Value = 0
$ws.Range("Y3").Formula = "=Y2+1"

$ws.Range("W4").Value = 275
$ws.Range("X4").Value = 0
$ws.Range("Y4").Formula = "=Y3+1"

$ws.Range("W5").Value = 319
$ws.Range("X5").Value = 0
$ws.Range("Y5").Formula = "=Y4+1"

$ws.Range("W6").Value = 269
$ws.Range("X6").Value = 0
$ws.Range("Y6").Formula = "=Y5+1"

# --- Rows 7:16 of V:Y held the old per-age-group rows -> clear them ------
$ws.Range("V7:Y16").ClearContents()

# --- Extend the blank W/X placeholder cells (rows 17-18) down to row 21,
#     formatted with the accounting/Comma (0-decimal) number format ------
$ws.Range("W17:X21").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# --- Selection / view state, matching the refreshed workbook -------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 16
$win.ScrollRow = 1
$ws.Range("W6").Select()

Write-Host "vce-youth.xlsx graph data refresh applied"
